{"js": "// 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\nconst dateResults = context.document.body.search(\"September 19, 2025\", { matchCase: true, matchWholeWord: false });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"September 21, 2025\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Split the recipient's mailing address line into two paragraphs:\n//    \"929 Story Road, San Jose CA 95122\" (non-table occurrence only)\n//      -> \"929 Story Road\"\n//      -> new paragraph \"San Jose, CA 95122\"\nconst addressResults = context.document.body.search(\"929 Story Road, San Jose CA 95122\", { matchCase: true, matchWholeWord: false });\naddressResults.load(\"items\");\nawait context.sync();\n\nconst addressParas = [];\nconst addressParentTables = [];\nfor (let i = 0; i < addressResults.items.length; i++) {\n  const para = addressResults.items[i].paragraphs.getFirst();\n  const parentTable = para.parentTableOrNullObject;\n  parentTable.load(\"isNullObject\");\n  addressParas.push(para);\n  addressParentTables.push(parentTable);\n}\nawait context.sync();\n\nlet addressIndex = -1;\nfor (let i = 0; i < addressResults.items.length; i++) {\n  if (addressParentTables[i].isNullObject) {\n    addressIndex = i;\n    break;\n  }\n}\n\nif (addressIndex !== -1) {\n  addressResults.items[addressIndex].insertText(\"929 Story Road\", Word.InsertLocation.replace);\n  addressParas[addressIndex].insertParagraph(\"San Jose, CA 95122\", Word.InsertLocation.after);\n}\nawait context.sync();\n\n// 3) Remove the now-redundant empty \"NoSpacing\" paragraph that follows\n//    the \"Board of Directors\" signature line.\nconst boardResults = context.document.body.search(\"Board of Directors\", { matchCase: true, matchWholeWord: false });\nboardResults.load(\"items\");\nawait context.sync();\n\nif (boardResults.items.length > 0) {\n  const boardPara = boardResults.items[0].paragraphs.getFirst();\n  const nextPara = boardPara.getNext();\n  nextPara.load(\"text,style\");\n  await context.sync();\n  if (nextPara.text === \"\" && nextPara.style === \"No Spacing\") {\n    nextPara.delete();\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"September 19, 2025\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"September 21, 2025\"\n$find.Execute($Missing, $Missing, $Missing, $Missing, $Missing, $Missing, $Missing, $Missing, $Missing, $Missing, 2)\n\n# 2) Split the recipient's mailing address line into two paragraphs.\n#    Only the body (non-table) occurrence of\n#    \"929 Story Road, San Jose CA 95122\" is affected; the identical text\n#    inside the \"PROPERTY ADDRESS\" table cell must stay untouched.\n$count = $d.Paragraphs.Count\n$addressPara = $null\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -eq \"929 Story Road, San Jose CA 95122`r\") {\n    if (-not $p.Range.Information(12)) {\n      $addressPara = $p\n      break\n    }\n  }\n}\nif ($addressPara -ne $null) {\n  $addressPara.Range.Text = \"929 Story Road`rSan Jose, CA 95122\"\n}\n\n# 3) Remove the now-redundant empty \"No Spacing\" paragraph that\n#    immediately follows the \"Board of Directors\" signature line.\n$count = $d.Paragraphs.Count\n$boardIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -eq \"Vietnam Town Condominium Owners Association Board of Directors`r\") {\n    $boardIndex = $i\n    break\n  }\n}\nif ($boardIndex -ge 1 -and $boardIndex -lt $d.Paragraphs.Count) {\n  $nextPara = $d.Paragraphs.Item($boardIndex + 1)\n  if ($nextPara.Range.Text -eq \"`r\" -and $nextPara.Style.NameLocal -eq \"No Spacing\") {\n    $nextPara.Range.Delete()\n  }\n}\n"}
